$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ntn1"
$ws.Cells.Item(2, 3).Value = "Mcam"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 2.082653666666667
$ws.Cells.Item(2, 8).Value = 6.247961
$ws.Cells.Item(2, 9).Value = 0.0472190032704503
$ws.Cells.Item(2, 10).Value = 0.0472190032704503
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 45.50099733333334
$ws.Cells.Item(2, 14).Value = 136.502992
$ws.Cells.Item(2, 15).Value = 0.5190633550775298
$ws.Cells.Item(2, 16).Value = 0.5190633550775298
$ws.Cells.Item(2, 17).Value = 94.76281893325691
$ws.Cells.Item(2, 18).Value = 852.865370399312
$ws.Cells.Item(2, 19).Value = 0.02450965426097678
$ws.Cells.Item(2, 20).Value = 0.02450965426097678

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ntn1"
$ws.Cells.Item(3, 3).Value = "Mcam"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 2.082653666666667
$ws.Cells.Item(3, 8).Value = 6.247961
$ws.Cells.Item(3, 9).Value = 0.0472190032704503
$ws.Cells.Item(3, 10).Value = 0.0472190032704503
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.5123886666666667
$ws.Cells.Item(3, 14).Value = 1.537166
$ws.Cells.Item(3, 15).Value = 0.005845194523436572
$ws.Cells.Item(3, 16).Value = 0.005845194523436572
$ws.Cells.Item(3, 17).Value = 1.067128135391778
$ws.Cells.Item(3, 18).Value = 9.604153218526001
$ws.Cells.Item(3, 19).Value = 0.0002760042593185697
$ws.Cells.Item(3, 20).Value = 0.0002760042593185697

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Ntn1"
$ws.Cells.Item(4, 3).Value = "Mcam"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 2.082653666666667
$ws.Cells.Item(4, 8).Value = 6.247961
$ws.Cells.Item(4, 9).Value = 0.0472190032704503
$ws.Cells.Item(4, 10).Value = 0.0472190032704503
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 41.64642833333333
$ws.Cells.Item(4, 14).Value = 124.939285
$ws.Cells.Item(4, 15).Value = 0.4750914503990336
$ws.Cells.Item(4, 16).Value = 0.4750914503990337
$ws.Cells.Item(4, 17).Value = 86.73508667198723
$ws.Cells.Item(4, 18).Value = 780.615780047885
$ws.Cells.Item(4, 19).Value = 0.02243334475015495
$ws.Cells.Item(4, 20).Value = 0.02243334475015495

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Ntn1"
$ws.Cells.Item(5, 3).Value = "Mcam"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 37.28222
$ws.Cells.Item(5, 8).Value = 111.84666
$ws.Cells.Item(5, 9).Value = 0.8452818134314446
$ws.Cells.Item(5, 10).Value = 0.8452818134314446
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 45.50099733333334
$ws.Cells.Item(5, 14).Value = 136.502992
$ws.Cells.Item(5, 15).Value = 0.5190633550775298
$ws.Cells.Item(5, 16).Value = 0.5190633550775298
$ws.Cells.Item(5, 17).Value = 1696.378192800747
$ws.Cells.Item(5, 18).Value = 15267.40373520672
$ws.Cells.Item(5, 19).Value = 0.4387548140657442
$ws.Cells.Item(5, 20).Value = 0.4387548140657442

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Ntn1"
$ws.Cells.Item(6, 3).Value = "Mcam"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 37.28222
$ws.Cells.Item(6, 8).Value = 111.84666
$ws.Cells.Item(6, 9).Value = 0.8452818134314446
$ws.Cells.Item(6, 10).Value = 0.8452818134314446
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.5123886666666667
$ws.Cells.Item(6, 14).Value = 1.537166
$ws.Cells.Item(6, 15).Value = 0.005845194523436572
$ws.Cells.Item(6, 16).Value = 0.005845194523436572
$ws.Cells.Item(6, 17).Value = 19.10298699617333
$ws.Cells.Item(6, 18).Value = 171.92688296556
$ws.Cells.Item(6, 19).Value = 0.004940836626630014
$ws.Cells.Item(6, 20).Value = 0.004940836626630014

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Ntn1"
$ws.Cells.Item(7, 3).Value = "Mcam"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 37.28222
$ws.Cells.Item(7, 8).Value = 111.84666
$ws.Cells.Item(7, 9).Value = 0.8452818134314446
$ws.Cells.Item(7, 10).Value = 0.8452818134314446
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 41.64642833333333
$ws.Cells.Item(7, 14).Value = 124.939285
$ws.Cells.Item(7, 15).Value = 0.4750914503990336
$ws.Cells.Item(7, 16).Value = 0.4750914503990337
$ws.Cells.Item(7, 17).Value = 1552.671303337567
$ws.Cells.Item(7, 18).Value = 13974.0417300381
$ws.Cells.Item(7, 19).Value = 0.4015861627390703
$ws.Cells.Item(7, 20).Value = 0.4015861627390704

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Ntn1"
$ws.Cells.Item(8, 3).Value = "Mcam"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 4.741387
$ws.Cells.Item(8, 8).Value = 14.224161
$ws.Cells.Item(8, 9).Value = 0.107499183298105
$ws.Cells.Item(8, 10).Value = 0.107499183298105
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 45.50099733333334
$ws.Cells.Item(8, 14).Value = 136.502992
$ws.Cells.Item(8, 15).Value = 0.5190633550775298
$ws.Cells.Item(8, 16).Value = 0.5190633550775298
$ws.Cells.Item(8, 17).Value = 215.7378372433014
$ws.Cells.Item(8, 18).Value = 1941.640535189712
$ws.Cells.Item(8, 19).Value = 0.05579888675080875
$ws.Cells.Item(8, 20).Value = 0.05579888675080875

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Ntn1"
$ws.Cells.Item(9, 3).Value = "Mcam"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 4.741387
$ws.Cells.Item(9, 8).Value = 14.224161
$ws.Cells.Item(9, 9).Value = 0.107499183298105
$ws.Cells.Item(9, 10).Value = 0.107499183298105
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.5123886666666667
$ws.Cells.Item(9, 14).Value = 1.537166
$ws.Cells.Item(9, 15).Value = 0.005845194523436572
$ws.Cells.Item(9, 16).Value = 0.005845194523436572
$ws.Cells.Item(9, 17).Value = 2.429432963080667
$ws.Cells.Item(9, 18).Value = 21.864896667726
$ws.Cells.Item(9, 19).Value = 0.0006283536374879878
$ws.Cells.Item(9, 20).Value = 0.0006283536374879876

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Ntn1"
$ws.Cells.Item(10, 3).Value = "Mcam"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 4.741387
$ws.Cells.Item(10, 8).Value = 14.224161
$ws.Cells.Item(10, 9).Value = 0.107499183298105
$ws.Cells.Item(10, 10).Value = 0.107499183298105
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 41.64642833333333
$ws.Cells.Item(10, 14).Value = 124.939285
$ws.Cells.Item(10, 15).Value = 0.4750914503990336
$ws.Cells.Item(10, 16).Value = 0.4750914503990337
$ws.Cells.Item(10, 17).Value = 197.4618338960983
$ws.Cells.Item(10, 18).Value = 1777.156505064885
$ws.Cells.Item(10, 19).Value = 0.05107194290980829
$ws.Cells.Item(10, 20).Value = 0.05107194290980829
